# Adds a new "2023" column (column I) to the sheet, mirroring the
# formatting of the existing "2022" column (column H), and removes the
# stale selection/active-cell marker from the saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column I, keyed by row number. $null marks rows whose
# H-column cell is present but empty (header/section rows) -- those still
# get a formatted-but-empty I cell to match column H.
$values = @{
    4  = 2023
    5  = 5.7627222366917641
    7  = 12.150374768642443
    8  = 0.33526865592135835
    9  = $null
    10 = 4.1862505557986136
    11 = 6.7483562655646434
    12 = $null
    13 = 11.721778533441505
    14 = 0.39226026012037718
    15 = 15.431347214780089
    16 = 12.56881331951053
    17 = 4.9361801817513591
    18 = 13.92788271827051
    19 = 5.796504268446359
    20 = 3.6469692666385813
    21 = 1.2344990530700553
    22 = $null
    23 = 0.33854574252686492
    24 = 7.0095823182535142
    25 = 7.4180588363268161
    26 = $null
    27 = 1.3575537444685963
    28 = 3.7000582818073822
    29 = 7.0145832826742662
    30 = 7.601388319014589
    31 = 3.2001807961995414
    32 = $null
    33 = 6.1374688939827911
    34 = 5.522716841454633
    35 = 5.959494359842247
    36 = 5.4831892692336535
    37 = 5.7612749525079918
}

foreach ($row in 4..37) {
    if (-not $values.ContainsKey($row)) { continue }

    $srcCell = $ws.Cells.Item($row, 8)   # column H
    $dstCell = $ws.Cells.Item($row, 9)   # column I

    # Copy column H's formatting (and value) into column I, then overwrite
    # with the real 2023 figure (or leave blank for section-header rows).
    $srcCell.Copy($dstCell)

    $newVal = $values[$row]
    if ($null -eq $newVal) {
        $dstCell.ClearContents()
    } else {
        $dstCell.Value = $newVal
    }
}

# The saved sheetView no longer pins the old D1 active-cell selection;
# reset it to the sheet's default top-left cell.
[void]$ws.Range("A1").Select()
